$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "表1" currently spans A1:R21 (20 data rows + header).
# We append 10 new rows (22-31): same parameter set as the existing
# "factor_apm_zscore_SD(0.0225)" rows, but for the new alpha
# "factor_apm_reg1d(wl=60,i=T)".
$newAlpha = "factor_apm_reg1d(wl=60,i=T)"

# Source rows to copy from (the existing factor_apm_* rows), in order.
$srcRows = 3,5,7,9,11,13,15,17,19,21
$destRow = 22

foreach ($srcRow in $srcRows) {
    $srcRange = $ws.Range("A" + $srcRow + ":R" + $srcRow)
    $destRange = $ws.Range("A" + $destRow + ":R" + $destRow)
    $srcRange.Copy($destRange)
    $ws.Range("B" + $destRow).Value = $newAlpha
    $destRow = $destRow + 1
}

$excel.CutCopyMode = 0

# Resize the table / autofilter to include the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:R31"))

# Update the active selection to mirror the post-edit UI state.
$ws.Range("A22").Select()
